$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 283, shifting existing rows 283:378 down to 285:380.
$ws.Rows("283:284").Insert()

# New row 283 (copy of the former row 283's static fields, with updated market data).
$ws.Cells.Item(283, 1).Value = 8
$ws.Cells.Item(283, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(283, 3).Value = "Coquimbo"
$ws.Cells.Item(283, 4).Value = 44900
$ws.Cells.Item(283, 5).Value = 4
$ws.Cells.Item(283, 6).Value = 100112003
$ws.Cells.Item(283, 7).Value = "Ajo"
$ws.Cells.Item(283, 8).Value = "Chino"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 400
$ws.Cells.Item(283, 11).Value = 16000
$ws.Cells.Item(283, 12).Value = 17000
$ws.Cells.Item(283, 13).Value = 16500
$ws.Cells.Item(283, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(283, 15).Value = "China"
$ws.Cells.Item(283, 16).Value = 1650
$ws.Cells.Item(283, 17).Value = 10
$ws.Cells.Item(283, 18).Value = "Hortaliza"

# New row 284 (copy of the former row 284's static fields, with updated market data).
$ws.Cells.Item(284, 1).Value = 8
$ws.Cells.Item(284, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(284, 3).Value = "Coquimbo"
$ws.Cells.Item(284, 4).Value = 44900
$ws.Cells.Item(284, 5).Value = 4
$ws.Cells.Item(284, 6).Value = 100112003
$ws.Cells.Item(284, 7).Value = "Ajo"
$ws.Cells.Item(284, 8).Value = "Chino"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 400
$ws.Cells.Item(284, 11).Value = 18000
$ws.Cells.Item(284, 12).Value = 19000
$ws.Cells.Item(284, 13).Value = 18500
$ws.Cells.Item(284, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(284, 15).Value = "China"
$ws.Cells.Item(284, 16).Value = 1850
$ws.Cells.Item(284, 17).Value = 10
$ws.Cells.Item(284, 18).Value = "Hortaliza"
